$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("signup")

# Update the registration/customer id column (D2:D9) from text codes (e.g. "123-1")
# to plain numeric credential values (123..130).
$values = @(123, 124, 125, 126, 127, 128, 129, 130)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Range("D$row").Value = $values[$i]
}

# Update the active selection to match the authored state (F6 instead of G6).
$ws.Range("F6").Select()
